$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Relocate the "total" row (old row 14) and the footer row (old row 15)
#    three rows down (to rows 17 and 18) to make room for three new product
#    rows. Merge destination ranges BEFORE pasting formats so the existing
#    cell styles are reused as-is instead of Excel re-deriving merge-aware
#    border variants.
# ---------------------------------------------------------------------------
$ws.Range("A18:E18").Merge()
$ws.Range("F18:G18").Merge()
$ws.Range("I18:N18").Merge()

$ws.Range("A15:N15").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A15:N15").Copy()
$ws.Range("A18:N18").PasteSpecial(-4163)   # xlPasteValues
$ws.Rows.Item(18).RowHeight = 16.5

$ws.Range("K17:N17").Merge()

$ws.Range("A14:N14").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A14:N14").Copy()
$ws.Range("A17:N17").PasteSpecial(-4163)   # xlPasteValues
$ws.Rows.Item(17).RowHeight = 26.25

# Clear the old merges that used to live on rows 14/15 - they are about to
# become ordinary product rows (14/15/16).
$ws.Range("K14:N14").UnMerge()
$ws.Range("A15:E15").UnMerge()
$ws.Range("F15:G15").UnMerge()
$ws.Range("I15:N15").UnMerge()

# ---------------------------------------------------------------------------
# 2) Build the three new product rows (14, 15, 16) by copying the formatting
#    of an existing product row (row 13 uses the shared style set 6/7/8/9)
# ---------------------------------------------------------------------------
$ws.Range("A13:N13").Copy()
$ws.Range("A14:N14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A13:N13").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A13:N13").Copy()
$ws.Range("A16:N16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()

$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

$ws.Range("B16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()

$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 3) Write the final product list (alphabetically ordered), rows 11-16.
#    A new syringe size, a new razor-blade item and a new cream item were
#    added; one pre-existing row's count also changed.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "سرنجات 10 سم"
$ws.Range("H11").Value = "-1:0"
$ws.Range("L11").Value = 4
$ws.Range("N11").Value = "1:0"

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "سرنجات 5 سم"
$ws.Range("H12").Value = "-1:0"
$ws.Range("L12").Value = 2
$ws.Range("N12").Value = "1:0"

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "شفرات فينوس حريمي "
$ws.Range("H13").Value = "16:0"
$ws.Range("L13").Value = 40
$ws.Range("N13").Value = "2:0"

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "كالونا "
$ws.Range("H14").Value = "-1:0"
$ws.Range("L14").Value = 15
$ws.Range("N14").Value = "1:0"

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "كريم فيرند لافلي الصغير"
$ws.Range("H15").Value = "6:0"
$ws.Range("L15").Value = 20
$ws.Range("N15").Value = "1:0"

$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "محلول ملح"
$ws.Range("H16").Value = "27:0"
$ws.Range("L16").Value = 48
$ws.Range("N16").Value = "2:0"

# ---------------------------------------------------------------------------
# 4) Update the grand total (moved from K14 to K17): sum of all L column
#    quantities (19+58+95+47+51+23+20+4+2+40+15+20+48 = 442)
# ---------------------------------------------------------------------------
$ws.Range("K17").Value = 442
